$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "26.158.19"
$ws.Range("E2").Value = "  +3.59%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.601.77"
$ws.Range("E3").Value = "  +2.87%  "

# Row 4
$ws.Range("E4").Value = "  -0.29%  "

# Row 5
Set-TextValue $ws.Range("D5") "212.88"
$ws.Range("E5").Value = "  +3.16%  "

# Row 6
$ws.Range("E6").Value = "  -0.32%  "

# Row 7
$ws.Range("E7").Value = "  +2.07%  "

# Row 8
$ws.Range("E8").Value = "  +3.50%  "

# Row 9
$ws.Range("E9").Value = "  +2.15%  "

# Row 10
Set-TextValue $ws.Range("D10") "18.01"
$ws.Range("E10").Value = "  +1.84%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0817"

# Row 12
Set-TextValue $ws.Range("D12") "1.824.86"
$ws.Range("E12").Value = "  +2.95%  "

# Row 13
Set-TextValue $ws.Range("D13") "1.604.11"
$ws.Range("E13").Value = "  +3.03%  "

# Row 14
$ws.Range("E14").Value = "  +0.53%  "

# Row 15
$ws.Range("E15").Value = "  +2.26%  "

# Row 16
Set-TextValue $ws.Range("D16") "26.154.06"

# Row 17
$ws.Range("E17").Value = "  +2.96%  "

# Row 18
$ws.Range("E18").Value = "  +2.08%  "

# Row 19
$ws.Range("E19").Value = "  -0.31%  "

# Row 20
Set-TextValue $ws.Range("D20") "205.14"
$ws.Range("E20").Value = "  +11.28%  "

# Row 21
$ws.Range("E21").Value = "  +3.81%  "

# Row 22
Set-TextValue $ws.Range("D22") "9.30"
$ws.Range("E22").Value = "  +0.95%  "

# Row 23
$ws.Range("E23").Value = "  +2.74%  "

# Row 24
Set-TextValue $ws.Range("D24") "1.83"
$ws.Range("E24").Value = "  +10.72%  "

# Row 25
Set-TextValue $ws.Range("D25") "141.98"
$ws.Range("E25").Value = "  +2.01%  "

# Row 26
$ws.Range("E26").Value = "  -0.33%  "

# Row 27
$ws.Range("E27").Value = "  -2.34%  "

# Row 28
Set-TextValue $ws.Range("D28") "15.22"
$ws.Range("E28").Value = "  +3.28%  "

# Row 29
Set-TextValue $ws.Range("D29") "6.44"
$ws.Range("E29").Value = "  +1.08%  "

# Row 30
$ws.Range("E30").Value = "  +1.88%  "

# Row 31
Set-TextValue $ws.Range("D31") "0.0471"
$ws.Range("E31").Value = "  +2.21%  "

# Row 32
$ws.Range("E32").Value = "  +4.22%  "

# Row 33
$ws.Range("E33").Value = "  +0.11%  "

# Row 34
$ws.Range("E34").Value = "  +2.35%  "

# Row 35
$ws.Range("E35").Value = "  +1.79%  "

# Row 36
Set-TextValue $ws.Range("D36") "1.112.16"
$ws.Range("E36").Value = "  +2.59%  "

# Row 37
$ws.Range("E37").Value = "  +9.10%  "

# Row 38
$ws.Range("E38").Value = "  -0.01%  "

# Row 39
Set-TextValue $ws.Range("D39") "2.32"
$ws.Range("E39").Value = "  +2.58%  "

# Row 40
$ws.Range("E40").Value = "  +2.79%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.493"
$ws.Range("E41").Value = "  +0.42%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.780"
$ws.Range("E42").Value = "  -4.08%  "

# Row 43
Set-TextValue $ws.Range("D43") "1.737.92"
$ws.Range("E43").Value = "  +2.99%  "

# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D44") "5.12"
$ws.Range("E44").Value = "  +1.72%  "

# Row 45
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D45") "92.86"
$ws.Range("E45").Value = "  +0.40%  "

# Row 46
$ws.Range("E46").Value = "  +5.23%  "

# Row 47
Set-TextValue $ws.Range("D47") "53.40"
$ws.Range("E47").Value = "  +2.35%  "

# Row 48
Set-TextValue $ws.Range("D48") "0.0503"
$ws.Range("E48").Value = "  -0.04%  "

# Row 49
$ws.Range("E49").Value = "  +0.93%  "

# Row 50
$ws.Range("E50").Value = "  -0.04%  "

# Row 51
Set-TextValue $ws.Range("D51") "7.20"
$ws.Range("E51").Value = "  +1.33%  "
